$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (443) down to the
# new rows (444-453) so the new rows match the existing table's look
# (text number format + borders for A/B, general number format + borders
# for C:F).
$ws.Range("A443:F443").Copy() | Out-Null
$ws.Range("A444:F453").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# New log rows for November 2024 (from saved_df_2024_11.csv)
$rows = @(
    @("2024-11-04 13:26:13", "020263.NC", 6, 19, 12, 1),
    @("2024-11-04 14:22:09", "020236.NC", 10, 50, 3, 1),
    @("2024-11-04 16:10:12", "020239.NC", 45, 31, 20, 1),
    @("2024-11-05 08:52:00", "020238.NC", 6, 7, 16, 1),
    @("2024-11-05 10:01:15", "020237.NC", 5, 18, 8, 1),
    @("2024-11-05 12:24:03", "L_15430.NC", 1, 5.4, 8, 1),
    @("2024-11-05 14:01:46", "L_15431.NC", 5, 7.5, 10, 1),
    @("2024-11-06 07:59:18", "L_15433.NC", 29, 143, 10, 1),
    @("2024-11-06 10:22:03", "L_15432.NC", 26, 91, 10, 1),
    @("2024-11-06 14:14:22", "020084.NC", 8, 1.6, 6, 1)
)

$startRow = 444
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

# Columns A and B no longer carry custom (wide) widths in the new layout;
# columns 1-6 all share the same width that columns C:F already had.
# (ColumnWidth is expressed in "characters"; column C's underlying raw
# width already serializes back out as the unmodified 8.85156 value, so we
# mirror its character width here rather than writing a literal to avoid
# introducing rounding drift.)
$stdWidth = $ws.Columns.Item(3).ColumnWidth
$ws.Columns.Item(1).ColumnWidth = $stdWidth
$ws.Columns.Item(2).ColumnWidth = $stdWidth
